# no-op baseline test
$wb = $excel.ActiveWorkbook
